$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1052990764"
$ws.Range("D16").Value = "GLORIA STEFANI ALVAREZ MAYA"
$ws.Range("E16").Value = "1607"
$ws.Range("F16").Value = 28859
$ws.Range("G16").Value = 781242
$ws.Range("C17").Value = "30579339"
$ws.Range("D17").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E17").Value = "1804"
$ws.Range("F17").Value = 32142
$ws.Range("G17").Value = 803535
$ws.Range("C18").Value = "30579339"
$ws.Range("D18").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E18").Value = "1805"
$ws.Range("F18").Value = 32142
$ws.Range("G18").Value = 803535
$ws.Range("C19").Value = "30579339"
$ws.Range("D19").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E19").Value = "1806"
$ws.Range("F19").Value = 32142
$ws.Range("G19").Value = 803535
$ws.Range("C20").Value = "30579339"
$ws.Range("D20").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E20").Value = "1807"
$ws.Range("F20").Value = 32142
$ws.Range("G20").Value = 803535
$ws.Range("C21").Value = "30579339"
$ws.Range("D21").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E21").Value = "1808"
$ws.Range("F21").Value = 32142
$ws.Range("G21").Value = 803535
$ws.Range("C22").Value = "30579339"
$ws.Range("D22").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E22").Value = "1809"
$ws.Range("F22").Value = 32142
$ws.Range("G22").Value = 803535
$ws.Range("C23").Value = "30579339"
$ws.Range("D23").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E23").Value = "1810"
$ws.Range("F23").Value = 32142
$ws.Range("G23").Value = 803535
$ws.Range("C24").Value = "30579339"
$ws.Range("D24").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E24").Value = "1811"
$ws.Range("F24").Value = 32142
$ws.Range("G24").Value = 803535
$ws.Range("C25").Value = "30579339"
$ws.Range("D25").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E25").Value = "1812"
$ws.Range("F25").Value = 32142
$ws.Range("G25").Value = 803535
$ws.Range("C26").Value = "30579339"
$ws.Range("D26").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E26").Value = "1901"
$ws.Range("F26").Value = 32142
$ws.Range("G26").Value = 803535
$ws.Range("C27").Value = "30579339"
$ws.Range("D27").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E27").Value = "1902"
$ws.Range("F27").Value = 32142
$ws.Range("G27").Value = 803535
$ws.Range("C28").Value = "30579339"
$ws.Range("D28").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E28").Value = "1903"
$ws.Range("F28").Value = 32142
$ws.Range("G28").Value = 803535
$ws.Range("C29").Value = "30579339"
$ws.Range("D29").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E29").Value = "1904"
$ws.Range("F29").Value = 32142
$ws.Range("G29").Value = 803535
$ws.Range("C30").Value = "30579339"
$ws.Range("D30").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E30").Value = "1905"
$ws.Range("F30").Value = 32142
$ws.Range("G30").Value = 803535
$ws.Range("C31").Value = "30579339"
$ws.Range("D31").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E31").Value = "1906"
$ws.Range("F31").Value = 32142
$ws.Range("G31").Value = 803535
$ws.Range("C32").Value = "30579339"
$ws.Range("D32").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E32").Value = "1907"
$ws.Range("F32").Value = 32142
$ws.Range("G32").Value = 803535
$ws.Range("C33").Value = "30579339"
$ws.Range("D33").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E33").Value = "1908"
$ws.Range("F33").Value = 32142
$ws.Range("G33").Value = 803535
$ws.Range("C34").Value = "30579339"
$ws.Range("D34").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E34").Value = "1909"
$ws.Range("F34").Value = 32142
$ws.Range("G34").Value = 803535
$ws.Range("C35").Value = "30579339"
$ws.Range("D35").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E35").Value = "1910"
$ws.Range("F35").Value = 32142
$ws.Range("G35").Value = 803535
$ws.Range("C36").Value = "30579339"
$ws.Range("D36").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E36").Value = "1911"
$ws.Range("F36").Value = 32142
$ws.Range("G36").Value = 803535
$ws.Range("C37").Value = "30579339"
$ws.Range("D37").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E37").Value = "1912"
$ws.Range("F37").Value = 32142
$ws.Range("G37").Value = 803535
$ws.Range("C38").Value = "30579339"
$ws.Range("D38").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E38").Value = "2001"
$ws.Range("F38").Value = 32142
$ws.Range("G38").Value = 803535
$ws.Range("C39").Value = "30579339"
$ws.Range("D39").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E39").Value = "2002"
$ws.Range("F39").Value = 32142
$ws.Range("G39").Value = 803535
$ws.Range("C40").Value = "30579339"
$ws.Range("D40").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E40").Value = "2003"
$ws.Range("F40").Value = 32142
$ws.Range("G40").Value = 803535
$ws.Range("C41").Value = "30579339"
$ws.Range("D41").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E41").Value = "2004"
$ws.Range("F41").Value = 32142
$ws.Range("G41").Value = 803535
$ws.Range("C42").Value = "30579339"
$ws.Range("D42").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E42").Value = "2005"
$ws.Range("F42").Value = 32142
$ws.Range("G42").Value = 803535
$ws.Range("C43").Value = "30579339"
$ws.Range("D43").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E43").Value = "2006"
$ws.Range("F43").Value = 32142
$ws.Range("G43").Value = 803535
$ws.Range("C44").Value = "30579339"
$ws.Range("D44").Value = "YANIS PATRICIA RODELO CATALAN"
$ws.Range("E44").Value = "2007"
$ws.Range("F44").Value = 32142
$ws.Range("G44").Value = 803535
